# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the most recently handed-back files (row 2, i.e. the
# 54b515d0-... entry) on both the "zh-cn" and "de-de" status sheets.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E2").Value = "2016-03-13 00:42:25"
$ws_zhcn.Range("H2").Value = "2016-03-13 00:42:41"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E2").Value = "2016-03-13 00:42:28"
$ws_dede.Range("H2").Value = "2016-03-13 00:42:47"
